# ============================================================================
# Restructure the three worksheets (SalesLog, Collections, Assignments) to
# match the updated "sales_ui" workbook layout:
#   - SalesLog:    drop StartDate/EndDate, add a single SentDate column,
#                  apply currency / percent / date number formats, drop one
#                  empty data row.
#   - Collections: replace DepositDue with a new CollectionDate column
#                  (inserted before Client), add currency / date formats and
#                  extra empty data rows.
#   - Assignments: move Client next to EndDate, append Notes / TaskStatus /
#                  Completed columns.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "SalesLog"
# Before: QuoteID,Client,QuotedPrice,Status,SalesRep,Deposit%,DepositPaid,
#         StartDate,EndDate,JobType
# After:  QuoteID,Client,QuotedPrice,Status,SalesRep,Deposit%,DepositPaid,
#         SentDate,JobType
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("SalesLog")

# Drop the last (now empty) data row - only 2 data rows remain afterwards.
$ws1.Rows.Item(4).Delete()

# Drop EndDate (column I); JobType shifts left from J into I.
$ws1.Columns.Item(9).Delete()

# StartDate (column H) becomes SentDate.
$ws1.Range("H1").Value = "SentDate"

# Number formats for the data rows.
$ws1.Range("C2:C3").NumberFormat = '"$"#,##0.00'
$ws1.Range("F2:F3").NumberFormat = '0.00"%"'
$ws1.Range("G2:G3").NumberFormat = '"$"#,##0.00'
$ws1.Range("H2:H3").NumberFormat = 'yyyy-mm-dd'

[void]$ws1.Range("A2:XFD4").Select()

# ----------------------------------------------------------------------
# Sheet "Collections"
# Before: QuoteID,Client,DepositDue,DepositPaid,BalanceDue,Status
# After:  QuoteID,CollectionDate,Client,DepositPaid,BalanceDue,Status
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Collections")

# Insert a new column at B for CollectionDate (Client etc. shift right).
$ws2.Columns.Item(2).Insert()
$ws2.Range("A1").Copy($ws2.Range("B1"))
$ws2.Range("B1").Value = "CollectionDate"

# DepositDue (old column C) is now column D - remove it entirely.
$ws2.Columns.Item(4).Delete()

# Add three empty data rows with the relevant number formats.
$ws2.Range("B2:B4").NumberFormat = 'yyyy-mm-dd'
$ws2.Range("D2:E4").NumberFormat = '"$"#,##0.00'

[void]$ws2.Range("A2:XFD6").Select()

# ----------------------------------------------------------------------
# Sheet "Assignments"
# Before: QuoteID,StartDate,EndDate,CrewMember,Payment,DaysTaken,Client
# After:  QuoteID,StartDate,EndDate,Client,CrewMember,Payment,DaysTaken,
#         Notes,TaskStatus,Completed
# ----------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Assignments")

# Insert a new column at D for Client (CrewMember.. shift right).
$ws3.Columns.Item(4).Insert()
$ws3.Range("A1").Copy($ws3.Range("D1"))
$ws3.Range("D1").Value = "Client"

# Old Client column has shifted from G to H - remove it.
$ws3.Columns.Item(8).Delete()

# Append the three new trailing columns.
$ws3.Range("A1").Copy($ws3.Range("H1"))
$ws3.Range("H1").Value = "Notes"
$ws3.Range("A1").Copy($ws3.Range("I1"))
$ws3.Range("I1").Value = "TaskStatus"
$ws3.Range("A1").Copy($ws3.Range("J1"))
$ws3.Range("J1").Value = "Completed"

# The insert pulled formatting into the new Client data cell - clear it so
# row 2 only keeps the StartDate/EndDate date-time formatting.
$ws3.Range("D2").Clear()
$ws3.Range("B2:C2").NumberFormat = 'yyyy\-mm\-dd\ hh:mm:ss'

[void]$ws3.Range("H16").Select()
